# Aggiunta caricamento di vari Sudoku da input, possibilita' dell'utente di
# vedere le possibili liste quando digita un numero sbagliato, controllo sui digits

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "controllo errore" (row 14) is now also assigned to the 2nd developer column (D) - "x"
$ws.Range("D14").Value = "x"

# New backlog row: "caricamento vari Sudoku", also assigned to column D ("x")
$ws.Range("C16").Value = "caricamento vari Sudoku"
$ws.Range("D16").Value = "x"

# Column C was widened to fit the longer new text
$ws.Columns.Item(3).ColumnWidth = 21.85

# Try to restore the window placement/size recorded for this edit (best effort;
# some runtimes may not persist these into the saved workbookView).
$win = $wb.Windows.Item(1)
$win.Left = 8700
$win.Top = 0
$win.Width = 10185
$win.Height = 10920
